# Insert a new weekly price record for "Vega Modelo de Temuco" / Puerro
# above the existing row 190. This shifts the old rows 190-193 down to
# 191-194 (preserving all of their original data), and the new row 190
# receives the new week's data (Fecha 2022-04-05 / serial 44656).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 190, pushing existing data down.
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with the new record.
$ws.Range("A190").Value = 10
$ws.Range("B190").Value = "Vega Modelo de Temuco"
$ws.Range("C190").Value = "La Araucanía"
$ws.Range("D190").Value = 44656
$ws.Range("E190").Value = 9
$ws.Range("F190").Value = 100112005
$ws.Range("G190").Value = "Puerro"
$ws.Range("H190").Value = "Azul de Maquehue"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 110
$ws.Range("K190").Value = 12000
$ws.Range("L190").Value = 12000
$ws.Range("M190").Value = 12000
$ws.Range("N190").Value = "`$/docena de paquetes"
$ws.Range("O190").Value = "Provincia de Cautín"
$ws.Range("P190").Value = 1000
$ws.Range("Q190").Value = 12
$ws.Range("R190").Value = "Hortaliza"
